$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 100; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2() -eq 45189) {
        $cell.Value = 45190
    }
}
